# Insert a new "submit" / "SUBMIT" key-value row above the existing
# "victory" row (row 12), pushing victory down to row 13, then move the
# active selection to the newly added B12 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "submit"
$ws.Range("B12").Value = "SUBMIT"

$ws.Range("B12").Select()
